# Update countries & provincias Spain
# Refresh of the COVID-19 "paises" data table:
#  - header timestamp cell updated
#  - several countries' stats refreshed
#  - Iran, Nepal and Malta each overtook the country(ies) right below them in the
#    ranking, so those rows now carry the overtaking country's fresh data while
#    the countries they passed keep their old (unchanged) numbers, shifted down
#    one row
#  - a brand-new row for "Islas Salomon" is appended at the bottom

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

# --- Estados Unidos (row 4) : updated totals -------------------------------
Set-Row 4 7549771 448 4777577 2558668 0 2 213526

# --- Reino Unido / Chile / Iran (rows 15-17): Iran overtakes both ----------
$ws.Cells.Item(15, 1).Value = "Iran"
Set-Row 15 468119 3523 387675 53698 0 179 26746

$ws.Cells.Item(16, 1).Value = "Reino Unido"
Set-Row 16 467146 0 0 0 0 0 42268

$ws.Cells.Item(17, 1).Value = "Chile"
Set-Row 17 466590 0 439607 14116 0 0 12867

# --- Rumania (row 32): updated totals --------------------------------------
Set-Row 32 134065 2064 107058 22060 0 32 4947

# --- Emiratos Arabes Unidos (row 44): updated totals ------------------------
Set-Row 44 97760 1231 87122 10212 0 2 426

# --- Japon / Nepal (rows 48-49): Nepal overtakes Japon ----------------------
$ws.Cells.Item(48, 1).Value = "Nepal"
Set-Row 48 84570 2120 62740 21302 0 8 528

$ws.Cells.Item(49, 1).Value = "Japon"
Set-Row 49 84215 0 77219 5418 0 0 1578

# --- Senegal (row 92): updated totals ---------------------------------------
Set-Row 92 15068 17 12751 2005 0 0 312

# --- Botsuana / Mali / Malta (rows 144-146): Malta overtakes both -----------
$ws.Cells.Item(144, 1).Value = "Malta"
Set-Row 144 3204 65 2711 455 0 1 38

$ws.Cells.Item(145, 1).Value = "Botsuana"
Set-Row 145 3172 0 710 2446 0 0 16

$ws.Cells.Item(146, 1).Value = "Mali"
Set-Row 146 3156 0 2467 558 0 0 131

# --- Islandia (row 148): updated totals -------------------------------------
Set-Row 148 2872 63 2212 650 0 0 10

# --- Gibraltar (row 180): updated totals ------------------------------------
Set-Row 180 428 12 358 70 0 0 0

# --- New country: Islas Salomon (row 220) -----------------------------------
$ws.Cells.Item(220, 1).Value = "Islas Salomon"
Set-Row 220 1 1 0 1 0 0 0

# --- Header timestamp (set last so the shared-string table appends the new
#     "Islas Salomon" entry before the title string gets its text refreshed) -
$ws.Range("A1").Value = "Datos actualizados a 3 de Octubre de 2020 a las 13:24"
